$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "Time (24HR)" header to "Time"
$ws.Range("D2").Value = "Time"

# Update the selected/active cell to D3 (matches the saved view state in the diff)
$ws.Range("D3").Select()
